# CIV-11262 Set sensitivity to public
#
# The authored change removes the "Classification: Controlled" content-
# marking text boxes that the Microsoft Purview / sensitivity-label add-in
# stamped into the even-page and first-page footers, and lowers the
# document's sensitivity label (tracked in docMetadata/LabelInfo.xml,
# which is maintained by the labeling service itself and isn't a
# scriptable part of the Word object model).
#
# Here we do the part that is reachable through the Word COM object
# model: delete the classification text-box shapes that were anchored in
# the footers, leaving the (now empty) footer paragraphs behind.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists -or $ftr.Shapes.Count -gt 0) {
            for ($j = $ftr.Shapes.Count; $j -ge 1; $j--) {
                $shp = $ftr.Shapes.Item($j)
                if ($shp.Name -like "Text Box*") {
                    $shp.Delete()
                }
            }
        }
    }
}
